$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "69.251.02"
Set-TextValue "E2" "  +0.92%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.776.08"
Set-TextValue "E3" "  -0.99%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.08%  "

# Row 5 - BNB
Set-TextValue "D5" "629.58"
Set-TextValue "E5" "  +3.38%  "

# Row 6 - Solana
Set-TextValue "D6" "165.97"
Set-TextValue "E6" "  +0.87%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.772.97"
Set-TextValue "E7" "  -1.02%  "

# Row 8 - USDC
Set-TextValue "E8" "  -0.02%  "

# Row 9 - XRP
Set-TextValue "D9" "0.522"
Set-TextValue "E9" "  +0.74%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.160"
Set-TextValue "E10" "  -0.37%  "

# Row 11 - Cardano
Set-TextValue "E11" "  +1.81%  "

# Row 12 - Toncoin
Set-TextValue "D12" "6.82"
Set-TextValue "E12" "  -2.61%  "

# Row 13 - ShibaInu
Set-TextValue "E13" "  -1.95%  "

# Row 14 - Avalanche
Set-TextValue "D14" "34.99"
Set-TextValue "E14" "  -0.72%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "4.410.55"
Set-TextValue "E15" "  -0.97%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.787.10"
Set-TextValue "E16" "  +0.15%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "69.254.85"
Set-TextValue "E17" "  +0.96%  "

# Row 18 - Chainlink
Set-TextValue "E18" "  -3.30%  "

# Row 19 - TRON
Set-TextValue "E19" "  -0.86%  "

# Row 20 - Polkadot
Set-TextValue "D20" "7.02"
Set-TextValue "E20" "  -0.91%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "464.36"
Set-TextValue "E21" "  +0.03%  "

# Row 22 - Uniswap
Set-TextValue "D22" "9.56"
Set-TextValue "E22" "  -0.76%  "

# Row 23 - Polygon
Set-TextValue "E23" "  +1.27%  "

# Row 24 - Litecoin
Set-TextValue "D24" "83.07"
Set-TextValue "E24" "  -0.84%  "

# Row 25 - PEPE
Set-TextValue "E25" "  -2.90%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue "D26" "11.98"
Set-TextValue "E26" "  -0.51%  "

# Row 27 - Fetch.AI
Set-TextValue "D27" "2.15"
Set-TextValue "E27" "  +1.48%  "

# Row 28 - was RenderToken, becomes Dai
Set-TextValue "B28" "Dai"
Set-TextValue "C28" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D28" "1.00"
Set-TextValue "E28" "  -0.01%  "

# Row 29 - was Dai, becomes RenderToken
Set-TextValue "B29" "RenderToken"
Set-TextValue "C29" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D29" "9.99"
Set-TextValue "E29" "  -0.51%  "

# Row 30 - WrappedeETH
Set-TextValue "D30" "3.925.04"
Set-TextValue "E30" "  -0.89%  "

# Row 31 - PancakeSwap
Set-TextValue "E31" "  +1.64%  "

# Row 32 - ImmutableX
Set-TextValue "D32" "2.26"
Set-TextValue "E32" "  +1.49%  "

# Row 33 - NEARProtocol
Set-TextValue "D33" "7.13"
Set-TextValue "E33" "  -2.17%  "

# Row 34 - EthereumClassic
Set-TextValue "D34" "28.62"
Set-TextValue "E34" "  -2.13%  "

# Row 35 - Kaspa
Set-TextValue "D35" "0.172"
Set-TextValue "E35" "  +15.67%  "

# Row 36 - Binance-PegBSC-USD
Set-TextValue "D36" "1.00"
Set-TextValue "E36" "  +0.09%  "

# Row 37 - was Aptos, becomes RenzoRestakedETH
Set-TextValue "B37" "RenzoRestakedETH"
Set-TextValue "C37" "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextValue "D37" "3.728.47"
Set-TextValue "E37" "  -0.88%  "

# Row 38 - was RenzoRestakedETH, becomes Aptos
Set-TextValue "B38" "Aptos"
Set-TextValue "C38" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D38" "8.98"
Set-TextValue "E38" "  -0.53%  "

# Row 39 - Hedera
Set-TextValue "E39" "  +0.34%  "

# Row 40 - dogwifhat
Set-TextValue "D40" "3.31"
Set-TextValue "E40" "  +1.88%  "

# Row 41 - Filecoin
Set-TextValue "D41" "5.80"
Set-TextValue "E41" "  -1.72%  "

# Row 42 - FirstDigitalUSD
Set-TextValue "E42" "  -0.09%  "

# Row 43 - Mantle
Set-TextValue "D43" "0.960"
Set-TextValue "E43" "  -2.11%  "

# Row 44 - USDe (unchanged)

# Row 45 - Monero
Set-TextValue "D45" "157.04"
Set-TextValue "E45" "  +2.54%  "

# Row 46 - ONDO
Set-TextValue "E46" "  +2.15%  "

# Row 47 - was Arweave, becomes Stacks
Set-TextValue "B47" "Stacks"
Set-TextValue "C47" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D47" "1.94"
Set-TextValue "E47" "  +3.65%  "

# Row 48 - was Stacks, becomes Arweave
Set-TextValue "B48" "Arweave"
Set-TextValue "C48" "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D48" "43.18"
Set-TextValue "E48" "  -0.19%  "

# Row 49 - TheGraph
Set-TextValue "D49" "0.295"
Set-TextValue "E49" "  -1.13%  "

# Row 50 - OKB
Set-TextValue "D50" "46.68"
Set-TextValue "E50" "  -1.03%  "

# Row 51 - Cosmos
Set-TextValue "E51" "  -0.37%  "
